# Update cryptos list with latest prices/volumes (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "66.251.97"
$ws.Cells.Item(2, 5).Value = "  +2.60%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.241.37"
$ws.Cells.Item(3, 5).Value = "  +5.13%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.02%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'574.34"
$ws.Cells.Item(5, 5).Value = "  +1.63%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'154.65"
$ws.Cells.Item(6, 5).Value = "  +8.30%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +0.05%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "3.238.05"
$ws.Cells.Item(8, 5).Value = "  +5.21%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.512"
$ws.Cells.Item(9, 5).Value = "  +3.53%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'7.16"
$ws.Cells.Item(10, 5).Value = "  +12.65%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.166"
$ws.Cells.Item(11, 5).Value = "  +4.63%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'0.484"
$ws.Cells.Item(12, 5).Value = "  +3.84%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'37.95"
$ws.Cells.Item(13, 5).Value = "  +6.03%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'0.0000234"
$ws.Cells.Item(14, 5).Value = "  +3.91%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.735.45"
$ws.Cells.Item(15, 5).Value = "  +4.59%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "66.165.29"
$ws.Cells.Item(16, 5).Value = "  +2.61%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'544.65"
$ws.Cells.Item(17, 5).Value = "  +9.79%  "

# Row 18
$ws.Cells.Item(18, 2).Value = "WrappedEther"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(18, 4).Value = "3.222.39"
$ws.Cells.Item(18, 5).Value = "  +4.60%  "

# Row 19
$ws.Cells.Item(19, 2).Value = "TRON"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(19, 4).Value = "'0.114"
$ws.Cells.Item(19, 5).Value = "  +2.71%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'7.03"
$ws.Cells.Item(20, 5).Value = "  +5.33%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'14.47"
$ws.Cells.Item(21, 5).Value = "  +4.95%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.735"
$ws.Cells.Item(22, 5).Value = "  +6.73%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'7.75"
$ws.Cells.Item(23, 5).Value = "  +6.94%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'13.46"
$ws.Cells.Item(24, 5).Value = "  +6.20%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'81.22"
$ws.Cells.Item(25, 5).Value = "  +3.49%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.28%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'9.40"
$ws.Cells.Item(27, 5).Value = "  +18.02%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'2.89"
$ws.Cells.Item(28, 5).Value = "  +3.99%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'2.26"
$ws.Cells.Item(29, 5).Value = "  +8.59%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'27.80"
$ws.Cells.Item(30, 5).Value = "  +4.80%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'2.79"
$ws.Cells.Item(31, 5).Value = "  +5.80%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -0.13%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +5.26%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'576.62"
$ws.Cells.Item(34, 5).Value = "  +11.73%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'5.78"
$ws.Cells.Item(35, 5).Value = "  +4.60%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'6.38"
$ws.Cells.Item(36, 5).Value = "  +7.17%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'0.0463"
$ws.Cells.Item(37, 5).Value = "  +13.98%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'54.37"
$ws.Cells.Item(38, 5).Value = "  +2.49%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Hedera"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(39, 4).Value = "'0.0864"
$ws.Cells.Item(39, 5).Value = "  +8.24%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "dogwifhat"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(40, 4).Value = "'3.08"
$ws.Cells.Item(40, 5).Value = "  +15.99%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  +4.43%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "3.145.90"
$ws.Cells.Item(42, 5).Value = "  +6.59%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'8.61"
$ws.Cells.Item(43, 5).Value = "  +2.93%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'2.36"
$ws.Cells.Item(44, 5).Value = "  +11.47%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.274"
$ws.Cells.Item(45, 5).Value = "  +11.26%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'26.69"
$ws.Cells.Item(46, 5).Value = "  +6.24%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +0.11%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "0.0₃0556"
$ws.Cells.Item(48, 5).Value = "  +2.41%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'0.113"
$ws.Cells.Item(49, 5).Value = "  +4.14%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'122.07"
$ws.Cells.Item(50, 5).Value = "  +0.37%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'2.24"
$ws.Cells.Item(51, 5).Value = "  +7.85%  "
